# "Fix Heatmap side pixels colorset"
# The B column held per-row distance/colour values used to build a side
# heatmap; the fix pins every row to the same computed pixel value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1:B100 all become the single fixed colorset pixel value.
$ws.Range("B1:B100").Value = 8660.25403784423

# The two "bestFit" columns (A = index, B = colorset value) re-fit to the
# new, now-uniform, string widths.
$ws.Columns.Item(1).ColumnWidth = 2.2916666666666665
$ws.Columns.Item(2).ColumnWidth = 10.916666666666666
